$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the next day's results as a new row at the bottom of the table.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

$ws.Cells.Item($newRow, 1).Value = 45980
$ws.Cells.Item($newRow, 2).Value = 68
$ws.Cells.Item($newRow, 3).Value = 77
$ws.Cells.Item($newRow, 4).Value = 77

# Match the date-cell style used by the rest of column A.
$ws.Cells.Item($newRow, 1).NumberFormat = $ws.Cells.Item($lastRow, 1).NumberFormat
